$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Profesor Asistente" row (old row 3) is merged into the "Profesor
# Asociado" row (row 2): the associate-professor stint now runs
# 2015-2019 instead of just 2019, so the separate 2015-2018 assistant
# row is removed and everything below shifts up one row.
$ws.Rows("3").Delete()

# Update the merged "Profesor Asociado" entry's date range and course.
$ws.Range("B2").Value2 = "2015 - 2019"

# Fix wording in the "Profesor Catedrático" course description.
$ws.Range("E6").Value2 = "Evolución y desarrollo de la comunicación vocal: cantos, modas y lenguaje (2016)"

# Update the merged course date range (kept for last so new shared
# strings are appended in the same order as the source edit).
$ws.Range("E2").Value2 = "Métodos cuantitativos II (Maestría en Psicología) (2017-2019)"

# Update the selected cell to match the saved selection in the workbook.
$ws.Range("E13").Select()
